# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
# Both sheets carry the same row-for-row event data, so the same set of
# F-column updates applies to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 746
    5  = 31
    7  = 229
    8  = 1682
    9  = 6371
    10 = 482
    11 = 364
    12 = 297
    13 = 95
    14 = 372
    15 = 137
    16 = 6328
    17 = 271
    18 = 1278
    20 = 117
    21 = 220
    22 = 104
    23 = 270
    24 = 104
    27 = 96
    28 = 8
    29 = 389
    30 = 87
    32 = 79
    35 = 23
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
